$wb = $excel.ActiveWorkbook
$count = $wb.Worksheets.Count
$last = $wb.Worksheets.Item($count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws.Name = "Aggregate"
$ws.Cells.Clear()

$ws.Cells.Item(1,1).Value = "FIPS"
$ws.Cells.Item(1,2).Value = "State"
$ws.Cells.Item(1,3).Value = "<12_month"
$ws.Cells.Item(1,4).Value = "At_risk"
$ws.Cells.Item(1,5).Value = "IR per 1000"

$ws.Cells.Item(2,1).Value = 4.0
$ws.Cells.Item(2,2).Value = "Arizona"
$ws.Cells.Item(2,3).Value = 42621.918836024626
$ws.Cells.Item(2,4).Value = 2802421.7470672023
$ws.Cells.Item(2,5).Value = 15.208959493919652

$ws.Cells.Item(3,1).Value = 6.0
$ws.Cells.Item(3,2).Value = "California"
$ws.Cells.Item(3,3).Value = 156599.0264181489
$ws.Cells.Item(3,4).Value = 16850453.394726697
$ws.Cells.Item(3,5).Value = 9.29346070101331

$ws.Cells.Item(4,1).Value = 9.0
$ws.Cells.Item(4,2).Value = "Connecticut"
$ws.Cells.Item(4,3).Value = 32938.5229378051
$ws.Cells.Item(4,4).Value = 2734477.7688722718
$ws.Cells.Item(4,5).Value = 12.045635664973537

$ws.Cells.Item(5,1).Value = 11.0
$ws.Cells.Item(5,2).Value = "District of Columbia"
$ws.Cells.Item(5,3).Value = 3184.446344618885
$ws.Cells.Item(5,4).Value = 179492.85599165817
$ws.Cells.Item(5,5).Value = 17.74135425627681

$ws.Cells.Item(6,1).Value = 13.0
$ws.Cells.Item(6,2).Value = "Geogia"
$ws.Cells.Item(6,3).Value = 94786.03809289802
$ws.Cells.Item(6,4).Value = 10458073.543508092
$ws.Cells.Item(6,5).Value = 9.063431969431596

$ws.Cells.Item(7,1).Value = 17.0
$ws.Cells.Item(7,2).Value = "Illinois"
$ws.Cells.Item(7,3).Value = 37799.24289699028
$ws.Cells.Item(7,4).Value = 5673571.265652148
$ws.Cells.Item(7,5).Value = 6.662336846956915

$ws.Cells.Item(8,1).Value = 18.0
$ws.Cells.Item(8,2).Value = "Indiana"
$ws.Cells.Item(8,3).Value = 105218.93189138013
$ws.Cells.Item(8,4).Value = 6936762.371679282
$ws.Cells.Item(8,5).Value = 15.168305652354107

$ws.Cells.Item(9,1).Value = 19.0
$ws.Cells.Item(9,2).Value = "Iowa"
$ws.Cells.Item(9,3).Value = 11510.486502884885
$ws.Cells.Item(9,4).Value = 1829733.9661702944
$ws.Cells.Item(9,5).Value = 6.290797851327419

$ws.Cells.Item(10,1).Value = 20.0
$ws.Cells.Item(10,2).Value = "Kansas"
$ws.Cells.Item(10,3).Value = 27509.487892604753
$ws.Cells.Item(10,4).Value = 3059760.3617240055
$ws.Cells.Item(10,5).Value = 8.990732815789757

$ws.Cells.Item(11,1).Value = 22.0
$ws.Cells.Item(11,2).Value = "Louisiana"
$ws.Cells.Item(11,3).Value = 5378.786540781787
$ws.Cells.Item(11,4).Value = 931966.3594570106
$ws.Cells.Item(11,5).Value = 5.771438514063552

$ws.Cells.Item(12,1).Value = 23.0
$ws.Cells.Item(12,2).Value = "Maine"
$ws.Cells.Item(12,3).Value = 6662.149701015229
$ws.Cells.Item(12,4).Value = 722763.1937653258
$ws.Cells.Item(12,5).Value = 9.217610634415294

$ws.Cells.Item(13,1).Value = 24.0
$ws.Cells.Item(13,2).Value = "Maryland"
$ws.Cells.Item(13,3).Value = 64870.5736086545
$ws.Cells.Item(13,4).Value = 5816583.904105316
$ws.Cells.Item(13,5).Value = 11.152692831073782

$ws.Cells.Item(14,1).Value = 26.0
$ws.Cells.Item(14,2).Value = "Michigan"
$ws.Cells.Item(14,3).Value = 126101.97053771566
$ws.Cells.Item(14,4).Value = 10491065.363787048
$ws.Cells.Item(14,5).Value = 12.0199394594369

$ws.Cells.Item(15,1).Value = 28.0
$ws.Cells.Item(15,2).Value = "Mississippi"
$ws.Cells.Item(15,3).Value = 18264.008952066375
$ws.Cells.Item(15,4).Value = 1300916.7417331804
$ws.Cells.Item(15,5).Value = 14.039337312035567

$ws.Cells.Item(16,1).Value = 29.0
$ws.Cells.Item(16,2).Value = "Missouri"
$ws.Cells.Item(16,3).Value = 46410.034572168624
$ws.Cells.Item(16,4).Value = 3600272.063188931
$ws.Cells.Item(16,5).Value = 12.890702079625912

$ws.Cells.Item(17,1).Value = 30.0
$ws.Cells.Item(17,2).Value = "Montana"
$ws.Cells.Item(17,3).Value = 3295.889765187216
$ws.Cells.Item(17,4).Value = 768012.0066596719
$ws.Cells.Item(17,5).Value = 4.291456040540417

$ws.Cells.Item(18,1).Value = 31.0
$ws.Cells.Item(18,2).Value = "Nebraska"
$ws.Cells.Item(18,3).Value = 18262.226573581298
$ws.Cells.Item(18,4).Value = 2014605.1364890952
$ws.Cells.Item(18,5).Value = 9.064916118206348

$ws.Cells.Item(19,1).Value = 33.0
$ws.Cells.Item(19,2).Value = "New Hampshire"
$ws.Cells.Item(19,3).Value = 9423.24961501254
$ws.Cells.Item(19,4).Value = 788301.5634397555
$ws.Cells.Item(19,5).Value = 11.953863917120971

$ws.Cells.Item(20,1).Value = 34.0
$ws.Cells.Item(20,2).Value = "New Jersey"
$ws.Cells.Item(20,3).Value = 51471.692726499445
$ws.Cells.Item(20,4).Value = 5274309.777003881
$ws.Cells.Item(20,5).Value = 9.758943805484705

$ws.Cells.Item(21,1).Value = 35.0
$ws.Cells.Item(21,2).Value = "New Mexico"
$ws.Cells.Item(21,3).Value = 8857.142946737338
$ws.Cells.Item(21,4).Value = 1327496.3453727677
$ws.Cells.Item(21,5).Value = 6.6720657858008705

$ws.Cells.Item(22,1).Value = 36.0
$ws.Cells.Item(22,2).Value = "New York"
$ws.Cells.Item(22,3).Value = 221226.36562273267
$ws.Cells.Item(22,4).Value = 15027480.545143578
$ws.Cells.Item(22,5).Value = 14.721454135851552

$ws.Cells.Item(23,1).Value = 39.0
$ws.Cells.Item(23,2).Value = "Ohio"
$ws.Cells.Item(23,3).Value = 71567.97481901593
$ws.Cells.Item(23,4).Value = 4755244.714276284
$ws.Cells.Item(23,5).Value = 15.050324246016872

$ws.Cells.Item(24,1).Value = 40.0
$ws.Cells.Item(24,2).Value = "Oklahoma"
$ws.Cells.Item(24,3).Value = 24627.721388718226
$ws.Cells.Item(24,4).Value = 2285659.3166476926
$ws.Cells.Item(24,5).Value = 10.774887232467767

$ws.Cells.Item(25,1).Value = 41.0
$ws.Cells.Item(25,2).Value = "Oregon"
$ws.Cells.Item(25,3).Value = 8328.024176806282
$ws.Cells.Item(25,4).Value = 752767.6262806169
$ws.Cells.Item(25,5).Value = 11.063207138641959

$ws.Cells.Item(26,1).Value = 42.0
$ws.Cells.Item(26,2).Value = "Pennsylvania"
$ws.Cells.Item(26,3).Value = 62291.57800717453
$ws.Cells.Item(26,4).Value = 4733924.715291146
$ws.Cells.Item(26,5).Value = 13.158548509646815

$ws.Cells.Item(27,1).Value = 44.0
$ws.Cells.Item(27,2).Value = "Rhode Island"
$ws.Cells.Item(27,3).Value = 5476.426426490039
$ws.Cells.Item(27,4).Value = 384116.9213563548
$ws.Cells.Item(27,5).Value = 14.257186085820528

$ws.Cells.Item(28,1).Value = 48.0
$ws.Cells.Item(28,2).Value = "Texas"
$ws.Cells.Item(28,3).Value = 381999.04549088026
$ws.Cells.Item(28,4).Value = 22992022.96837785
$ws.Cells.Item(28,5).Value = 16.61441648767765

$ws.Cells.Item(29,1).Value = 49.0
$ws.Cells.Item(29,2).Value = "Utah"
$ws.Cells.Item(29,3).Value = 30221.406332381404
$ws.Cells.Item(29,4).Value = 2902955.358223645
$ws.Cells.Item(29,5).Value = 10.410565304343592

$ws.Cells.Item(30,1).Value = 50.0
$ws.Cells.Item(30,2).Value = "Vermont"
$ws.Cells.Item(30,3).Value = 6498.410650110485
$ws.Cells.Item(30,4).Value = 563279.5559912089
$ws.Cells.Item(30,5).Value = 11.536741536225586

$ws.Cells.Item(31,1).Value = 53.0
$ws.Cells.Item(31,2).Value = "Washington"
$ws.Cells.Item(31,3).Value = 18647.25128699706
$ws.Cells.Item(31,4).Value = 2752372.813938709
$ws.Cells.Item(31,5).Value = 6.774972922477175

$ws.Cells.Item(32,1).Value = 54.0
$ws.Cells.Item(32,2).Value = "West Virginia"
$ws.Cells.Item(32,3).Value = 3846.644645155366
$ws.Cells.Item(32,4).Value = 325031.2108367552
$ws.Cells.Item(32,5).Value = 11.834693152244132

$ws.Cells.Item(33,1).Value = 55.0
$ws.Cells.Item(33,2).Value = "Wisconsin"
$ws.Cells.Item(33,3).Value = 14403.531708608236
$ws.Cells.Item(33,4).Value = 1174446.719861135
$ws.Cells.Item(33,5).Value = 12.264099737373606
